$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe in the literal forces Excel to store the value as
# text (quotePrefix) instead of auto-converting it to a number, matching
# the source data where these price/percent columns are stored as text.

$ws.Range("D2").Value = '''61.077.92'
$ws.Range("E2").Value = '''  +0.89%  '
$ws.Range("D3").Value = '''3.334.51'
$ws.Range("E3").Value = '''  +0.38%  '
$ws.Range("E4").Value = '''  -0.07%  '
$ws.Range("D5").Value = '''400.81'
$ws.Range("E5").Value = '''  -2.20%  '
$ws.Range("D6").Value = '''126.17'
$ws.Range("E6").Value = '''  +9.81%  '
$ws.Range("D7").Value = '''0.591'
$ws.Range("E7").Value = '''  +3.35%  '
$ws.Range("D8").Value = '''0.998'
$ws.Range("E8").Value = '''  -0.18%  '
$ws.Range("D9").Value = '''0.659'
$ws.Range("E9").Value = '''  +5.55%  '
$ws.Range("E10").Value = '''  +3.47%  '
$ws.Range("D11").Value = '''41.05'
$ws.Range("E11").Value = '''  +2.84%  '
$ws.Range("E12").Value = '''  -0.75%  '
$ws.Range("D13").Value = '''3.863.53'
$ws.Range("E13").Value = '''  +0.49%  '
$ws.Range("E14").Value = '''  +1.44%  '
$ws.Range("D15").Value = '''19.30'
$ws.Range("E15").Value = '''  +1.03%  '
$ws.Range("D16").Value = '''3.308.66'
$ws.Range("E16").Value = '''  -0.59%  '
$ws.Range("D17").Value = '''60.893.95'
$ws.Range("E17").Value = '''  +0.73%  '
$ws.Range("D18").Value = '''11.24'
$ws.Range("E18").Value = '''  +4.43%  '
$ws.Range("E19").Value = '''  +0.72%  '
$ws.Range("E20").Value = '''  +11.63%  '
$ws.Range("E21").Value = '''  -4.89%  '
$ws.Range("D22").Value = '''80.30'
$ws.Range("E22").Value = '''  +8.62%  '
$ws.Range("D23").Value = '''12.86'
$ws.Range("E23").Value = '''  +3.70%  '
$ws.Range("D24").Value = '''299.88'
$ws.Range("E24").Value = '''  +1.55%  '
$ws.Range("D25").Value = '''3.10'
$ws.Range("E25").Value = '''  -0.33%  '
$ws.Range("D26").Value = '''4.67'
$ws.Range("E26").Value = '''  +10.19%  '
$ws.Range("D27").Value = '''8.30'
$ws.Range("E27").Value = '''  +10.77%  '
$ws.Range("D28").Value = '''28.96'
$ws.Range("E28").Value = '''  -0.72%  '
$ws.Range("E29").Value = '''  -2.81%  '
$ws.Range("E30").Value = '''  -0.18%  '
$ws.Range("E31").Value = '''  -0.17%  '
$ws.Range("D32").Value = '''11.40'
$ws.Range("E32").Value = '''  +1.34%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '''  -0.04%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").Value = '''2.53'
$ws.Range("E34").Value = '''  +2.73%  '
$ws.Range("D35").Value = '''41.01'
$ws.Range("E35").Value = '''  +0.31%  '
$ws.Range("D36").Value = '''0.0475'
$ws.Range("E36").Value = '''  -2.91%  '
$ws.Range("D37").Value = '''52.12'
$ws.Range("E37").Value = '''  +0.03%  '
$ws.Range("D38").Value = '''0.996'
$ws.Range("E38").Value = '''  -0.21%  '
$ws.Range("D39").Value = '''3.37'
$ws.Range("E39").Value = '''  +0.57%  '
$ws.Range("E40").Value = '''  -4.79%  '
$ws.Range("E41").Value = '''  +4.37%  '
$ws.Range("D42").Value = '''136.05'
$ws.Range("E42").Value = '''  +1.56%  '
$ws.Range("E43").Value = '''  +2.27%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.280'
$ws.Range("E44").Value = '''  -3.68%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''3.88'
$ws.Range("E45").Value = '''  +1.24%  '
$ws.Range("D46").Value = '''16.67'
$ws.Range("E46").Value = '''  +2.75%  '
$ws.Range("E47").Value = '''  +0.41%  '
$ws.Range("D48").Value = '''21.24'
$ws.Range("E48").Value = '''  +0.94%  '
$ws.Range("D49").Value = '''2.119.00'
$ws.Range("E49").Value = '''  -0.74%  '
$ws.Range("D50").Value = '''3.664.31'
$ws.Range("E50").Value = '''  +0.67%  '
$ws.Range("E51").Value = '''  -1.52%  '
